$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FFMP rows (155-158) - fill in Question/Answer (columns A/B) first
$ws.Range("A155").Value2 = "FFMP"
$ws.Range("B155").Value2 = "Flexible Flow Management Plan"

$ws.Range("A156").Value2 = "FFMP"
$ws.Range("B156").Value2 = "Address: Water Supply, Instream Flow Needs, and Spill Migration"

$ws.Range("A157").Value2 = "FFMP"
$ws.Range("B157").Value2 = "Based on Reservoir Releases"

$ws.Range("A158").Value2 = "FFMP"
$ws.Range("B158").Value2 = "Renewal of FFMP requires unanimous vote from PA, NJ, NY, DE, and NYC"

$ws.Range("A155:B158").WrapText = $true

# Backfill the "Flow Management" category (column C) for all the rows in
# this section (110-158), matching the newly-introduced Category column.
$catRange = $ws.Range("C110:C158")
$catRange.Value2 = "Flow Management"
$catRange.WrapText = $true

# New Non-drought Flow Targets rows (159-160) - no Category column
$ws.Range("A159").Value2 = "Non-drought Flow Targets"
$ws.Range("B159").Value2 = "Montague - 1750 cfs (1130 MGD)"

$ws.Range("A160").Value2 = "Non-drought Flow Targets"
$ws.Range("B160").Value2 = "Trenton - 3000 cfs (1940 MGD)"

$ws.Range("A159:B160").WrapText = $true

# Row heights to mirror the multi-line wrapped text sizing used elsewhere
# in the sheet (17 ~ 1 line, 34 ~ 2 lines).
$ws.Rows.Item(155).RowHeight = 17
$ws.Rows.Item(156).RowHeight = 34
$ws.Rows.Item(157).RowHeight = 17
$ws.Rows.Item(158).RowHeight = 34
$ws.Rows.Item(159).RowHeight = 17
$ws.Rows.Item(160).RowHeight = 17

# Update the view selection to match the post-edit state.
$ws.Range("D161").Select()
